# Updates cryptocurrency price and 1h volume-change figures in the
# "cryptos" worksheet, mirroring the upstream GitHub Actions scraper
# refresh (commit: "Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.068.42"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "3.508.23"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.92"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.34"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("E9").Value = "  +4.78%  "
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.434"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "4.116.73"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.09"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").Value = "67.054.99"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").Value = "3.491.57"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.33"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.17"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.84"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.05"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.13"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.537"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.22"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.33"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.85"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.38"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.63"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.884"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.91"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.09"
$ws.Range("E39").Value = "  +4.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.69"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0748"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.42"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.46"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "2.814.04"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.58"
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.73"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "340.71"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.74"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.49"
$ws.Range("E51").Value = "  -0.79%  "
